$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight row 9 (A9:I9) with a solid yellow fill, matching the new style (fillId=2)
$ws.Range("A9:I9").Interior.Color = 65535

# Update row 10 values (iteration changed 50000 -> 10000, learning_rate 0.3 -> 0.5, recomputed losses)
$ws.Range("B10").Value = 10000
$ws.Range("C10").Value = 0.5
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 128
$ws.Range("H10").Value = 0.064
$ws.Range("I10").Value = 0.16

# New row 11
$ws.Range("B11").Value = 100000
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 128
$ws.Range("H11").Value = 0.047
$ws.Range("I11").Value = 0.154

# New row 12
$ws.Range("B12").Value = 50000
$ws.Range("C12").Value = 0.3
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 128
$ws.Range("H12").Value = 0.05
$ws.Range("I12").Value = 0.19

# New row 13
$ws.Range("B13").Value = 50000
$ws.Range("C13").Value = 0.3
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 128
$ws.Range("H13").Value = 0.048
$ws.Range("I13").Value = 0.154

# New row 14
$ws.Range("B14").Value = 100000
$ws.Range("C14").Value = 0.1
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 128
$ws.Range("H14").Value = 0.049
$ws.Range("I14").Value = 0.219

# Update selection to match target (D17)
$ws.Range("D17").Select()
